$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates as described by the diff.
# For numeric-looking "Price" strings (column D) that Excel would
# otherwise auto-convert to a floating point number (losing the
# original text formatting / introducing binary rounding), force
# the cell to Text first, write the value, then restore the style
# to Normal so no stray number-format style is left behind.

$ws.Range('D2').Value = '69.834.04'
$ws.Range('E2').Value = '  -1.24%  '
$ws.Range('D3').Value = '3.725.54'
$ws.Range('E3').Value = '  -2.35%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '616.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').Value = '3.723.21'
$ws.Range('E7').Value = '  -2.34%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.531'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.165'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.32'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.32%  '
$ws.Range('E12').Value = '  -5.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.03'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.18%  '
$ws.Range('E14').Value = '  -4.31%  '
$ws.Range('D15').Value = '4.329.00'
$ws.Range('E15').Value = '  -2.66%  '
$ws.Range('D16').Value = '3.708.04'
$ws.Range('E16').Value = '  -2.47%  '
$ws.Range('D17').Value = '69.800.83'
$ws.Range('E17').Value = '  -1.47%  '
$ws.Range('E18').Value = '  -2.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.01%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '501.20'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.56%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.28'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.31'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.27%  '
$ws.Range('E23').Value = '  -3.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.52'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.39'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.98'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.17'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.06%  '
$ws.Range('E28').Value = '  +3.86%  '
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.46'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.24%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.91'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.97'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '30.36'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.115'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.08'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.139'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.343'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.07'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.06'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '49.99'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '426.67'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.93'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.27%  '
$ws.Range('E45').Value = '  -3.82%  '
$ws.Range('D46').Value = '2.947.32'
$ws.Range('E46').Value = '  -6.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0360'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '27.20'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '136.21'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.46'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.63%  '
